$wb = $excel.ActiveWorkbook

# 1. Rename "Hoja3" (the waterfall/cascada sheet) to "Cascada"
$cascada = $wb.Worksheets.Item("Hoja3")
$cascada.Name = "Cascada"

# 2. Add a new sheet after "Cascada" and name it "Hoja1"
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $cascada)
$newSheet.Name = "Hoja1"

# 3. Populate the data table used by the radar chart
$newSheet.Range("B1").Value = 2014
$newSheet.Range("C1").Value = 2010
$newSheet.Range("D1").Value = 2006

$newSheet.Range("A2").Value = "Activos"
$newSheet.Range("A3").Value = "Ocupados"
$newSheet.Range("A4").Value = "Inactivos"
$newSheet.Range("A5").Value = "Parados"

$newSheet.Range("B2").Value = 22954.6
$newSheet.Range("C2").Value = 23364.6
$newSheet.Range("D2").Value = 21780

$newSheet.Range("B3").Value = 17344.2
$newSheet.Range("C3").Value = 18724.5
$newSheet.Range("D3").Value = 19939.1

$newSheet.Range("B4").Value = 15560
$newSheet.Range("C4").Value = 15395.5
$newSheet.Range("D4").Value = 15362.9

$newSheet.Range("B5").Value = 5610.4
$newSheet.Range("C5").Value = 4640.1000000000004
$newSheet.Range("D5").Value = 1840.9

$newSheet.Range("B2:D6").NumberFormat = "#,##0.00"

Write-Output "done basic data"
